# Update Name of Algo - apply corrected values to result_data_KNN sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 13.33
$ws.Range("E4").Value = 12.919
$ws.Range("E7").Value = 13.35
$ws.Range("E8").Value = 12.913
$ws.Range("A11").Value = -21.76
$ws.Range("A12").Value = -21.776
$ws.Range("E12").Value = 13.143
$ws.Range("E14").Value = 12.911
$ws.Range("A15").Value = -21.279
$ws.Range("E22").Value = 12.862
